$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '61.404.13'
$ws.Range('E2').Value = '  +0.81%  '

# Row 3
$ws.Range('D3').Value = '2.932.45'
$ws.Range('E3').Value = '  +0.66%  '

# Row 4
$ws.Range('E4').Value = '  -0.08%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.60'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.21%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.83'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.20%  '

# Row 7
$ws.Range('E7').Value = '  +0.00%  '

# Row 8
$ws.Range('E8').Value = '  -1.23%  '

# Row 9
$ws.Range('E9').Value = '  +0.59%  '

# Row 10
$ws.Range('E10').Value = '  -1.30%  '

# Row 11
$ws.Range('E11').Value = '  -0.45%  '

# Row 12
$ws.Range('E12').Value = '  -0.79%  '

# Row 13
$ws.Range('E13').Value = '  -0.87%  '

# Row 14
$ws.Range('E14').Value = '  +0.50%  '

# Row 15
$ws.Range('D15').Value = '3.418.77'
$ws.Range('E15').Value = '  +0.57%  '

# Row 16
$ws.Range('D16').Value = '61.394.53'
$ws.Range('E16').Value = '  +0.79%  '

# Row 17
$ws.Range('D17').Value = '2.933.38'
$ws.Range('E17').Value = '  +0.57%  '

# Row 18
$ws.Range('E18').Value = '  -0.76%  '

# Row 19
$ws.Range('E19').Value = '  +0.87%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.55'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.48%  '

# Row 21
$ws.Range('E21').Value = '  -0.86%  '

# Row 22
$ws.Range('E22').Value = '  +0.20%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '81.36'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.14%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.84'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.00%  '

# Row 25
$ws.Range('E25').Value = '  -1.12%  '

# Row 26
$ws.Range('E26').Value = '  -1.15%  '

# Row 27
$ws.Range('E27').Value = '  -0.03%  '

# Row 28
$ws.Range('E28').Value = '  -3.61%  '

# Row 29
$ws.Range('E29').Value = '  -0.53%  '

# Row 30
$ws.Range('E30').Value = '  -2.17%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '26.75'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.04%  '

# Row 32
$ws.Range('E32').Value = '  +1.51%  '

# Row 33
$ws.Range('E33').Value = '  -0.02%  '

# Row 34
$ws.Range('E34').Value = '  +2.49%  '

# Row 35
$ws.Range('E35').Value = '  -0.05%  '

# Row 36
$ws.Range('E36').Value = '  -0.13%  '

# Row 37
$ws.Range('E37').Value = '  -1.00%  '

# Row 38
$ws.Range('E38').Value = '  +0.50%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.122'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.04%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.48'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.56%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '41.90'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.76%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.279'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.34%  '

# Row 43
$ws.Range('E43').Value = '  +0.01%  '

# Row 44
$ws.Range('D44').Value = '2.695.22'
$ws.Range('E44').Value = '  -0.21%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '133.38'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.76%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '363.55'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.90%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.54'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.96%  '

# Row 49
$ws.Range('E49').Value = '  -1.32%  '

# Row 50
$ws.Range('E50').Value = '  -0.48%  '

# Row 51
$ws.Range('E51').Value = '  +0.47%  '
